$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) dates for rows 2-9 from 45208 to 45212
$ws.Range("C2:C9").Value = 45212
